# Apply scheduled-runner price/profit updates to the Gungnir Profits workbook.
# Each block targets one job-sheet and rewrites the currentAveragePrice* / Leve Price*
# / LeveProfit* columns (H:N) for the rows the runner refreshed.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 866.55554
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 999.8570999999999
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 2999.5713
$ws.Range("M121").Value = 547
$ws.Range("N121").Value = -6493.5713
$ws.Range("H137").Value = 2566.0833
$ws.Range("I137").Value = 2784.7144
$ws.Range("J137").Value = 2260
$ws.Range("K137").Value = 8354.143199999999
$ws.Range("L137").Value = 6780
$ws.Range("M137").Value = -5804.143199999999
$ws.Range("N137").Value = -11880
$ws.Range("H138").Value = 2873.9285
$ws.Range("I138").Value = 2340.1667
$ws.Range("J138").Value = 3274.25
$ws.Range("K138").Value = 7020.500100000001
$ws.Range("L138").Value = 9822.75
$ws.Range("M138").Value = -1880.500100000001
$ws.Range("N138").Value = -20102.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 168276.33
$ws.Range("I45").Value = 334182.66
$ws.Range("J45").Value = 2370
$ws.Range("K45").Value = 334182.66
$ws.Range("L45").Value = 2370
$ws.Range("M45").Value = -333805.66
$ws.Range("N45").Value = -3124
$ws.Range("H82").Value = 13800
$ws.Range("J82").Value = 13800
$ws.Range("L82").Value = 13800
$ws.Range("N82").Value = -14522
$ws.Range("H85").Value = 13800
$ws.Range("J85").Value = 13800
$ws.Range("L85").Value = 13800
$ws.Range("N85").Value = -16296
$ws.Range("H110").Value = 1465.4286
$ws.Range("I110").Value = 901.8261
$ws.Range("J110").Value = 4058
$ws.Range("K110").Value = 901.8261
$ws.Range("L110").Value = 4058
$ws.Range("M110").Value = 1143.1739
$ws.Range("N110").Value = -8148
$ws.Range("H132").Value = 1032731.1
$ws.Range("I132").Value = 652.6875
$ws.Range("J132").Value = 6537149.5
$ws.Range("K132").Value = 1958.0625
$ws.Range("L132").Value = 19611448.5
$ws.Range("M132").Value = 571.9375
$ws.Range("N132").Value = -19616508.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 71430370
$ws.Range("I105").Value = 1798.091
$ws.Range("J105").Value = 333335100
$ws.Range("K105").Value = 1798.091
$ws.Range("L105").Value = 333335100
$ws.Range("M105").Value = -51.09099999999989
$ws.Range("N105").Value = -333338594
$ws.Range("H134").Value = 2648879.8
$ws.Range("I134").Value = 901.15625
$ws.Range("J134").Value = 11122411
$ws.Range("K134").Value = 2703.46875
$ws.Range("L134").Value = 33367233
$ws.Range("M134").Value = -168.46875
$ws.Range("N134").Value = -33372303

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1793260.1
$ws.Range("I31").Value = 2778643.2
$ws.Range("J31").Value = 1654.4546
$ws.Range("K31").Value = 2778643.2
$ws.Range("L31").Value = 1654.4546
$ws.Range("M31").Value = -2778348.2
$ws.Range("N31").Value = -2244.4546
$ws.Range("H34").Value = 1793260.1
$ws.Range("I34").Value = 2778643.2
$ws.Range("J34").Value = 1654.4546
$ws.Range("K34").Value = 2778643.2
$ws.Range("L34").Value = 1654.4546
$ws.Range("M34").Value = -2778441.2
$ws.Range("N34").Value = -2058.4546
$ws.Range("H134").Value = 1495.0588
$ws.Range("I134").Value = 1301.2307
$ws.Range("J134").Value = 2125
$ws.Range("K134").Value = 3903.6921
$ws.Range("L134").Value = 6375
$ws.Range("M134").Value = -1368.6921
$ws.Range("N134").Value = -11445

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 103.57143
$ws.Range("I2").Value = 105.2
$ws.Range("J2").Value = 99.5
$ws.Range("K2").Value = 631.2
$ws.Range("L2").Value = 597
$ws.Range("M2").Value = -518.2
$ws.Range("N2").Value = -823
$ws.Range("H131").Value = 820.29
$ws.Range("J131").Value = 820.49493
$ws.Range("L131").Value = 2461.48479
$ws.Range("N131").Value = -12541.48479
$ws.Range("H132").Value = 22733178
$ws.Range("I132").Value = 881.875
$ws.Range("J132").Value = 35723064
$ws.Range("K132").Value = 7936.875
$ws.Range("L132").Value = 321507576
$ws.Range("M132").Value = -5406.875
$ws.Range("N132").Value = -321512636
$ws.Range("H137").Value = 12196427
$ws.Range("I137").Value = 27778438
$ws.Range("J137").Value = 1810.1305
$ws.Range("K137").Value = 83335314
$ws.Range("L137").Value = 5430.3915
$ws.Range("M137").Value = -83330214
$ws.Range("N137").Value = -15630.3915
$ws.Range("H140").Value = 10418206
$ws.Range("I140").Value = 13158997
$ws.Range("K140").Value = 39476991
$ws.Range("M140").Value = -39471811

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 15376.6
$ws.Range("J57").Value = 15376.6
$ws.Range("L57").Value = 15376.6
$ws.Range("N57").Value = -17016.6
$ws.Range("H102").Value = 3628.4285
$ws.Range("I102").Value = 1549.5
$ws.Range("J102").Value = 4460
$ws.Range("K102").Value = 1549.5
$ws.Range("L102").Value = 4460
$ws.Range("M102").Value = 72.5
$ws.Range("N102").Value = -7704
$ws.Range("H113").Value = 1220
$ws.Range("I113").Value = 1220
$ws.Range("K113").Value = 1220
$ws.Range("M113").Value = 950
$ws.Range("H122").Value = 55568130
$ws.Range("I122").Value = 100021070
$ws.Range("J122").Value = 1952
$ws.Range("K122").Value = 300063210
$ws.Range("L122").Value = 5856
$ws.Range("M122").Value = -300060760
$ws.Range("N122").Value = -10756
$ws.Range("H132").Value = 5722.5713
$ws.Range("I132").Value = 2040.95
$ws.Range("J132").Value = 14926.625
$ws.Range("K132").Value = 6122.85
$ws.Range("L132").Value = 44779.875
$ws.Range("M132").Value = -3592.85
$ws.Range("N132").Value = -49839.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5845
$ws.Range("I22").Value = 5250
$ws.Range("K22").Value = 5250
$ws.Range("M22").Value = -4955
$ws.Range("H27").Value = 5845
$ws.Range("I27").Value = 5250
$ws.Range("K27").Value = 5250
$ws.Range("M27").Value = -5143
$ws.Range("H36").Value = 59775
$ws.Range("J36").Value = 59775
$ws.Range("L36").Value = 59775
$ws.Range("N36").Value = -60899
$ws.Range("H40").Value = 125002500
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 125002500
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -125002772
$ws.Range("H69").Value = 37053
$ws.Range("J69").Value = 37053
$ws.Range("L69").Value = 37053
$ws.Range("N69").Value = -38675
$ws.Range("H72").Value = 37053
$ws.Range("J72").Value = 37053
$ws.Range("L72").Value = 111159
$ws.Range("N72").Value = -119271
$ws.Range("H109").Value = 20129.5
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H122").Value = 7210.4736
$ws.Range("I122").Value = 7388.8335
$ws.Range("K122").Value = 22166.5005
$ws.Range("M122").Value = -19716.5005
$ws.Range("H136").Value = 194808510
$ws.Range("I136").Value = 114289360
$ws.Range("J136").Value = 1000000000
$ws.Range("K136").Value = 342868080
$ws.Range("L136").Value = 3000000000
$ws.Range("M136").Value = -342865530
$ws.Range("N136").Value = -3000005100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 33790.938
$ws.Range("I122").Value = 57127.777
$ws.Range("J122").Value = 3786.4285
$ws.Range("K122").Value = 171383.331
$ws.Range("L122").Value = 11359.2855
$ws.Range("M122").Value = -168933.331
$ws.Range("N122").Value = -16259.2855
$ws.Range("H126").Value = 1452.6
$ws.Range("I126").Value = 724.36365
$ws.Range("K126").Value = 2173.09095
$ws.Range("M126").Value = 296.9090500000002
$ws.Range("H136").Value = 2111.878
$ws.Range("I136").Value = 1250.6
$ws.Range("J136").Value = 2932.1428
$ws.Range("K136").Value = 3751.8
$ws.Range("L136").Value = 8796.428400000001
$ws.Range("M136").Value = -1201.8
$ws.Range("N136").Value = -13896.4284
